$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B3 to the numeric value 3434 (replacing the previous string "hghghg")
$ws.Range("B3").Value = 3434

# Update the selected cell to C4
$ws.Range("C4").Select()
